$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Perbaikan import data pegawai ---

# 1) Tanggal lahir default hasil import salah ketik: 2022-07-21 -> 2022-07-22.
#    String ini dipakai bersama (shared string) oleh kolom F "Tgl. Lahir" di
#    hampir semua baris pegawai, jadi perbaiki di setiap sel yang memakainya.
#    NumberFormat dipaksa ke teks dulu supaya nilainya tidak otomatis dibaca
#    sebagai tanggal oleh Excel, lalu format dikembalikan (ClearFormats) agar
#    selnya tetap tanpa style khusus seperti aslinya.
$rows = @(4) + @(6..190)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.NumberFormat = "@"
    $cell.Value = "2022-07-22"
    $cell.ClearFormats()
}

# 2) Baris 191 (pegawai "Burhanu Sultan Ramadan", NIP 2012.22.07.178) adalah
#    hasil import yang keliru/duplikat, hapus barisnya.
$ws.Rows.Item(191).Delete()
